# Customer.xlsx edit script
# Applies:
#  - model sheet: insert a new "verbose_name" column (D), shifting the
#    existing "through"..."forms" columns one to the right, fill it with
#    the Django verbose_name values, and convert the old boolean (True/False)
#    marker cells into quote-prefixed text cells ("True"/"False") in their
#    shifted positions.
#  - model_functions sheet: simplify the Meta class snippet text and shrink
#    the row that holds it.
#  - selection/view bookkeeping to match where the author left the cursor.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "model_functions" sheet -- edited first so the in-place rewrite of
#    the existing "Meta class" shared string keeps its original index;
#    new strings introduced below (verbose_name values, True/False) are
#    then appended after it, matching the author's original edit order.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("model_functions")

$ws2.Range("B3").Value = "    class Meta:`n        verbose_name = 'Customer'`n        verbose_name_plural = 'Customers'"
$ws2.Rows("3:3").RowHeight = 51

$ws2.Range("B10").Select()

# ---------------------------------------------------------------------
# 2. "model" sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("model")

# Insert a new column before D; this shifts old D..W to E..X and widens
# the dimension/used range accordingly.
$ws.Columns("D:D").Insert()
$ws.Range("D1").ColumnWidth = 31.33

# Header for the newly inserted column.
$ws.Range("D1").Value = "verbose_name"

# verbose_name values, one per data row (quote characters are part of the
# literal Python source text being modeled).
$ws.Range("D2").Value = "'" + '"Company Name"'
$ws.Range("D3").Value = '"Active"'
$ws.Range("D4").Value = '"AFM"'
$ws.Range("D5").Value = '"First Name"'
$ws.Range("D6").Value = '"Last Name"'
$ws.Range("D7").Value = '"Email"'
$ws.Range("D8").Value = '"Phone"'
$ws.Range("D9").Value = '"Address"'
$ws.Range("D10").Value = '"Created at"'
$ws.Range("D11").Value = '"Updated at"'

# Re-key the numeric max_length markers that were in column K (now L).
$ws.Range("L2").Value = 100
$ws.Range("L4").Value = 9
$ws.Range("L5").Value = 100
$ws.Range("L6").Value = 100
$ws.Range("L8").Value = 15

# The old boolean (TRUE/FALSE) marker cells, re-entered as quote-prefixed
# text in their shifted columns (old column + 1 letter).
$ws.Range("T2").Value = "'True"
$ws.Range("U2").Value = "'True"
$ws.Range("V2").Value = "'True"
$ws.Range("X2").Value = "'True"

$ws.Range("P3").Value = "'False"
$ws.Range("T3").Value = "'True"
$ws.Range("V3").Value = "'True"
$ws.Range("X3").Value = "'True"

$ws.Range("N4").Value = "'True"
$ws.Range("O4").Value = "'True"
$ws.Range("T4").Value = "'True"
$ws.Range("V4").Value = "'True"
$ws.Range("X4").Value = "'True"

$ws.Range("T5").Value = "'True"
$ws.Range("V5").Value = "'True"
$ws.Range("X5").Value = "'True"

$ws.Range("T6").Value = "'True"
$ws.Range("V6").Value = "'True"
$ws.Range("X6").Value = "'True"

$ws.Range("M7").Value = "'True"
$ws.Range("T7").Value = "'True"
$ws.Range("V7").Value = "'True"
$ws.Range("X7").Value = "'True"

$ws.Range("N8").Value = "'True"
$ws.Range("O8").Value = "'True"
$ws.Range("T8").Value = "'True"
$ws.Range("V8").Value = "'True"
$ws.Range("X8").Value = "'True"

$ws.Range("N9").Value = "'True"
$ws.Range("O9").Value = "'True"
$ws.Range("T9").Value = "'True"
$ws.Range("V9").Value = "'True"
$ws.Range("X9").Value = "'True"

$ws.Range("Q10").Value = "'True"
$ws.Range("T10").Value = "'True"
$ws.Range("V10").Value = "'True"
$ws.Range("X10").Value = "'True"
$ws.Range("X10").ClearContents()

$ws.Range("R11").Value = "'True"
$ws.Range("T11").Value = "'True"
$ws.Range("V11").Value = "'True"
$ws.Range("X11").Value = "'True"
$ws.Range("X11").ClearContents()

# Cursor position left by the author.
$ws.Range("L4").Select()
